# ------------------------------------------------------------------
# Updated Team Diary and worked on Wrangling
#
# 1) Fix a typo ("Al" -> "All") in the "Members present" column for the
#    10/21/2023 meeting row.
# 2) Correct the bogus Start/End time values on the 10/19 and 10/21 rows
#    (they previously held a full date-time serial instead of a plain
#    elapsed-time fraction).
# 3) Fill in the previously-blank row with the 10/23 meeting, and append
#    two further meetings (10/25 and 10/26) documenting continued work
#    on wrangling / temporal analysis.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) fix "Al" -> "All" typo -------------------------------------
$ws.Range("D9").Value = "All"

# --- 2) correct the Start/End time values already on the sheet -----
$ws.Range("B8").Value = 1.75
$ws.Range("C8").Value = 0
$ws.Range("B9").Value = 0.58333333333333337
$ws.Range("C9").Value = 1.8055555555555554

# --- 3) fill the blank row 10 with the 10/23 meeting ----------------
$ws.Range("A10").Value = 45222
$ws.Range("A10").NumberFormat = "m/d/yyyy"
$ws.Range("B10").Value = 0.66666666666666663
$ws.Range("C10").Value = 0.83333333333333337
$ws.Range("B10:C10").NumberFormat = "h:mm"
$ws.Range("D10").Value = "All"
$ws.Range("E10").Value = "Attempted to Solve Wrangling Problems"

# give row 10 the same look (font/fill/border/alignment) as the row above it
$ws.Range("A9:E9").Copy()
$ws.Range("A10:E10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
# ... but keep the date (not text) + right number formats from above
$ws.Range("A10").NumberFormat = "m/d/yyyy"
$ws.Range("B10:C10").NumberFormat = "h:mm"
$ws.Rows.Item(10).RowHeight = 16.95

# --- new rows: 10/25 and 10/26 meetings -----------------------------
$ws.Range("A11").Value = 45224
$ws.Range("B11").Value = 0.875
$ws.Range("C11").Value = 0.125
$ws.Range("D11").Value = "All"
$ws.Range("E11").Value = "Attempted to Solve Wrangling Problems + Attempted Temporal Analysis"

$ws.Range("A12").Value = 45225
$ws.Range("B12").Value = 0.64583333333333337
$ws.Range("C12").Value = 0.91666666666666663
$ws.Range("D12").Value = "All"
$ws.Range("E12").Value = "Discussed fatal problems regarding wrangling and attempted to solve them "

# match formatting of rows 11 & 12 to the rest of the table
$ws.Range("A10:E10").Copy()
$ws.Range("A11:E12").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A11:A12").NumberFormat = "m/d/yyyy"
$ws.Range("B11:C12").NumberFormat = "h:mm"
$ws.Rows.Item(11).RowHeight = 15.75
$ws.Rows.Item(12).RowHeight = 15.75

# --- cosmetic touch-ups to match the table's thin borders everywhere,
#     wrap text on the long "Discussions" column and column widths ----
# re-color the thin borders to a normal black instead of the odd
# indexed color the original (Google-Sheets-exported) file had
$ws.Range("A1:E12").Borders.LineStyle = 1   # xlContinuous
$ws.Range("A1:E12").Borders.Weight = 2      # xlThin
$ws.Range("A1:E12").Borders.Color = 0       # RGB black
$ws.Range("E7:E12").WrapText = $true

$ws.Columns.Item(2).ColumnWidth = 15.6
$ws.Columns.Item(3).ColumnWidth = 14.6
$ws.Columns.Item(5).ColumnWidth = 65.6

# leave the selection where the author left it when they saved the file
$ws.Range("G15").Select()
